$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "column C" values: drop the accidentally-concatenated
# "Bestand 1" data that had been glued onto the end of each cell, keeping
# only the (small) piece that genuinely belongs to "bestand 2", plus in a
# few rows the single matching "Bestand 1 column N" value. ---
$ws.Range("C1").Value = "Data bestand 2"
$ws.Range("C2").Value = "dwa"
$ws.Range("C3").Value = "dwadwa"
$ws.Range("C4").Value = "dwadwaBestand 1 column 4"
$ws.Range("C5").Value = "dada"
$ws.Range("C6").Value = "dadaBestand 1 column 6"
$ws.Range("C7").Value = "adaBestand 1 column 2"
$ws.Range("C8").Value = "ada"
$ws.Range("C9").Value = "dadaBestand 1 column 8"
$ws.Range("C10").Value = "dadaBestand 1 column 11"
$ws.Range("C11").Value = "dadaBestand 1 column 10"
$ws.Range("C12").Value = "adaBestand 1 column 9"

# --- Formatting feedback: cells that no longer need to stand out lose the
# highlight fill entirely, the rest keep being highlighted but in the new
# blue color (RGB 0, 105, 255 => 0x0069FF) instead of the old dark red. ---
$ws.Range("C1").ClearFormats()
$ws.Range("C2").ClearFormats()
$ws.Range("C3").ClearFormats()
$ws.Range("C5").ClearFormats()
$ws.Range("C8").ClearFormats()

$newHighlight = 0 + (105 * 256) + (255 * 65536)
$ws.Range("C4").Interior.Color = $newHighlight
$ws.Range("C6").Interior.Color = $newHighlight
$ws.Range("C7").Interior.Color = $newHighlight
$ws.Range("C9").Interior.Color = $newHighlight
$ws.Range("C10").Interior.Color = $newHighlight
$ws.Range("C11").Interior.Color = $newHighlight
$ws.Range("C12").Interior.Color = $newHighlight
